# Slide 9 ("Distribution") - add a footnote textbox below the glucose chart
# explaining the clinical benchmark values for Glucose.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# Position/size below come from EMU values (727650, 4552622, 7074373, 276999)
# converted to points (1 pt = 12700 EMU), since Shapes.AddTextbox works in points.
$box = $s.Shapes.AddTextbox(1, 57.2952755905512, 358.474173228346, 557.037244094488, 21.8109448818898)
$box.Name = "TextBox 1"
$box.Fill.Visible = 0

$tf = $box.TextFrame
$tf.WordWrap = 0
$tf.AutoSize = 1

$tr = $tf.TextRange
$tr.Text = "*For reference, a Glucose value of >140 is considered abnormal and a value >200 indicates diabetes"
$tr.Font.Size = 12

# Split into 3 runs matching the authored edit (middle run covers "of >140 ")
$midRun = $tr.Characters(33, 8)
$midRun.Text = "of >140 "
$midRun.Font.Size = 12
